# Made Experimenter and Comment Metrics objects
$wb = $excel.ActiveWorkbook

# --- Experimenter sheet: add data_reference, linked_references, description ---
$wsExperimenter = $wb.Worksheets.Item("Experimenter")
$wsExperimenter.Range("C1").Value = "data_reference"
$wsExperimenter.Range("D1").Value = "linked_references"
$wsExperimenter.Range("E1").Value = "description"

# --- Comment sheet: add data_reference, linked_references, name, description ---
$wsComment = $wb.Worksheets.Item("Comment")
$wsComment.Range("E1").Value = "data_reference"
$wsComment.Range("F1").Value = "linked_references"
$wsComment.Range("G1").Value = "name"
$wsComment.Range("H1").Value = "description"
